$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell D4: "Copy of #1 but ran 24 hours" -> "Copy of #1 but ran 48 hours"
$ws.Range("D4").Value = "Copy of #1 but ran 48 hours"

# Apply the existing date style (from B2) to the new date cells, then fill in values
$ws.Range("B2").Copy()
$ws.Range("B5:B7").PasteSpecial(-4122) # xlPasteFormats

# Row 5 (Test 3)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "2/18/2022"
$ws.Range("C5").Value = "sevy"
$ws.Range("D5").Value = "copy of #1 but with additional library with 10 new thermos"
$ws.Range("E5").Value = "/work/westgroup/harris.se/autoscience/nhept_iter1"

# Row 6 (Test 4)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "5/2/2022"
$ws.Range("C6").Value = "sevy"
$ws.Range("E6").Value = "/work/westgroup/harris.se/autoscience/nheptane4"
$ws.Range("D6").Value = "copy of #1 but using latest RMG-Py/RMG-database, and 7 day time limit"

# Row 7 (Test 5)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "5/2/2022"
$ws.Range("C7").Value = "sevy"
$ws.Range("D7").Value = "copy of #1 but using latest RMG-Py/RMG-database, and 24 hour time limit"
$ws.Range("E7").Value = "/work/westgroup/harris.se/autoscience/nheptane5"

$ws.Range("E8").Select()

$wb.Save()
